$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column (D) keeps its original text representation instead of
# being auto-converted to a number by Excel when values look numeric.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '51.927.02'
$ws.Range('E2').Value = '  +0.46%  '

$ws.Range('D3').Value = '2.935.07'
$ws.Range('E3').Value = '  +3.88%  '

$ws.Range('E4').Value = '  +0.05%  '

$ws.Range('D5').Value = '352.27'
$ws.Range('E5').Value = '  +0.67%  '

$ws.Range('D6').Value = '112.40'
$ws.Range('E6').Value = '  -0.08%  '

$ws.Range('D7').Value = '0.561'
$ws.Range('E7').Value = '  +0.86%  '

$ws.Range('E8').Value = '  +0.02%  '

$ws.Range('D9').Value = '0.627'
$ws.Range('E9').Value = '  +1.66%  '

$ws.Range('D10').Value = '39.38'
$ws.Range('E10').Value = '  -1.65%  '

$ws.Range('D11').Value = '0.0890'
$ws.Range('E11').Value = '  +5.12%  '

$ws.Range('E12').Value = '  +1.21%  '

$ws.Range('D13').Value = '19.99'
$ws.Range('E13').Value = '  +0.51%  '

$ws.Range('D14').Value = '7.83'
$ws.Range('E14').Value = '  +1.07%  '

$ws.Range('D15').Value = '3.398.11'
$ws.Range('E15').Value = '  +3.99%  '

$ws.Range('D16').Value = '2.924.21'
$ws.Range('E16').Value = '  +3.51%  '

$ws.Range('D17').Value = '0.989'
$ws.Range('E17').Value = '  +0.98%  '

$ws.Range('D18').Value = '51.993.81'
$ws.Range('E18').Value = '  +0.52%  '

$ws.Range('D19').Value = '7.64'
$ws.Range('E19').Value = '  +0.81%  '

$ws.Range('D20').Value = '3.31'
$ws.Range('E20').Value = '  -3.58%  '

$ws.Range('D21').Value = '14.28'
$ws.Range('E21').Value = '  +6.85%  '

$ws.Range('D22').Value = '0.0₃0986'
$ws.Range('E22').Value = '  +1.54%  '

$ws.Range('D23').Value = '71.22'

$ws.Range('D24').Value = '269.12'
$ws.Range('E24').Value = '  +0.30%  '

$ws.Range('D25').Value = '2.78'
$ws.Range('E25').Value = '  +1.34%  '

$ws.Range('E26').Value = '  +10.00%  '

$ws.Range('D27').Value = '26.96'
$ws.Range('E27').Value = '  +2.97%  '

$ws.Range('E28').Value = '  +0.12%  '

$ws.Range('D29').Value = '7.42'
$ws.Range('E29').Value = '  +16.96%  '

$ws.Range('E30').Value = '  +21.07%  '

$ws.Range('E31').Value = '  +0.42%  '

$ws.Range('D32').Value = '37.45'
$ws.Range('E32').Value = '  -2.48%  '

$ws.Range('E33').Value = '  +0.13%  '

$ws.Range('E34').Value = '  +10.90%  '

$ws.Range('D35').Value = '52.88'

$ws.Range('D36').Value = '0.0453'
$ws.Range('E36').Value = '  +1.25%  '

$ws.Range('D37').Value = '0.999'
$ws.Range('E37').Value = '  -0.10%  '

$ws.Range('D38').Value = '3.31'
$ws.Range('E38').Value = '  +3.13%  '

$ws.Range('E39').Value = '  +0.17%  '

$ws.Range('E40').Value = '  +2.11%  '

$ws.Range('D41').Value = '2.70'
$ws.Range('E41').Value = '  +7.66%  '

$ws.Range('E42').Value = '  +1.77%  '

$ws.Range('E43').Value = '  +5.96%  '

$ws.Range('E44').Value = '  -0.91%  '

$ws.Range('E45').Value = '  +1.09%  '

$ws.Range('E46').Value = '  +1.44%  '

$ws.Range('D47').Value = '2.170.09'
$ws.Range('E47').Value = '  +0.08%  '

$ws.Range('D48').Value = '111.76'
$ws.Range('E48').Value = '  -8.93%  '

$ws.Range('E49').Value = '  +0.63%  '

$ws.Range('D50').Value = '0.0344'
$ws.Range('E50').Value = '  +11.29%  '

$ws.Range('D51').Value = '0.938'
$ws.Range('E51').Value = '  -1.03%  '
